# Adds 42 new transaction rows (rows 314-355) to the 'Konto' sheet for the
# week ending 2021-03-14 (dates 2021-03-08 .. 2021-03-14), following the
# existing Datum/Receipt Number/Konto/Beskrivning/Debet/Kredit layout.
# - Column A: date serials formatted like the existing rows (copy A313's format).
# - Column B: Receipt Number; kept as text (leading apostrophe) so purely
#   numeric receipt numbers are not auto-converted to the Number type,
#   matching how the source data stores them. Blank cells use a lone
#   apostrophe to produce an explicit empty text value instead of Excel
#   clearing the cell outright.
# - Column C: Konto code (number), blank via apostrophe on the two rows
#   without a resolved account.
# - Column D: description text.
# - Column E/F: Debet/Kredit numeric amounts, blank (apostrophe) on the
#   side not used for a given line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = $ws.Range("A313").NumberFormat

# Row 314
$ws.Range("A314").NumberFormat = $dateFormat
$ws.Range("A314").Value = 44263
$ws.Range("B314").Value = "'1081133"
$ws.Range("C314").Value = 3011
$ws.Range("D314").Value = "Order 1081133 Swish +46734869113"
$ws.Range("E314").Value = "'"
$ws.Range("F314").Value = 396.43

# Row 315
$ws.Range("A315").NumberFormat = $dateFormat
$ws.Range("A315").Value = 44263
$ws.Range("B315").Value = "'1081133"
$ws.Range("C315").Value = 2611
$ws.Range("D315").Value = "Order 1081133 Swish +46734869113"
$ws.Range("E315").Value = "'"
$ws.Range("F315").Value = 47.57

# Row 316
$ws.Range("A316").NumberFormat = $dateFormat
$ws.Range("A316").Value = 44263
$ws.Range("B316").Value = "'1081133"
$ws.Range("C316").Value = 1930
$ws.Range("D316").Value = "Order 1081133 Swish +46734869113"
$ws.Range("E316").Value = 444
$ws.Range("F316").Value = "'"

# Row 317
$ws.Range("A317").NumberFormat = $dateFormat
$ws.Range("A317").Value = 44263
$ws.Range("B317").Value = "Reko19"
$ws.Range("C317").Value = 3011
$ws.Range("D317").Value = "Reko Swish +46739012974"
$ws.Range("E317").Value = "'"
$ws.Range("F317").Value = 282.14

# Row 318
$ws.Range("A318").NumberFormat = $dateFormat
$ws.Range("A318").Value = 44263
$ws.Range("B318").Value = "Reko19"
$ws.Range("C318").Value = 2611
$ws.Range("D318").Value = "Reko Swish +46739012974"
$ws.Range("E318").Value = "'"
$ws.Range("F318").Value = 33.86

# Row 319
$ws.Range("A319").NumberFormat = $dateFormat
$ws.Range("A319").Value = 44263
$ws.Range("B319").Value = "Reko19"
$ws.Range("C319").Value = 1930
$ws.Range("D319").Value = "Reko Swish +46739012974"
$ws.Range("E319").Value = 316
$ws.Range("F319").Value = "'"

# Row 320
$ws.Range("A320").NumberFormat = $dateFormat
$ws.Range("A320").Value = 44263
$ws.Range("B320").Value = "Reko20"
$ws.Range("C320").Value = 3011
$ws.Range("D320").Value = "Reko Swish +46704972332"
$ws.Range("E320").Value = "'"
$ws.Range("F320").Value = 141.07

# Row 321
$ws.Range("A321").NumberFormat = $dateFormat
$ws.Range("A321").Value = 44263
$ws.Range("B321").Value = "Reko20"
$ws.Range("C321").Value = 2611
$ws.Range("D321").Value = "Reko Swish +46704972332"
$ws.Range("E321").Value = "'"
$ws.Range("F321").Value = 16.93

# Row 322
$ws.Range("A322").NumberFormat = $dateFormat
$ws.Range("A322").Value = 44263
$ws.Range("B322").Value = "Reko20"
$ws.Range("C322").Value = 1930
$ws.Range("D322").Value = "Reko Swish +46704972332"
$ws.Range("E322").Value = 158
$ws.Range("F322").Value = "'"

# Row 323
$ws.Range("A323").NumberFormat = $dateFormat
$ws.Range("A323").Value = 44264
$ws.Range("B323").Value = "Reko21"
$ws.Range("C323").Value = 3011
$ws.Range("D323").Value = "Reko Swish +46709927597"
$ws.Range("E323").Value = "'"
$ws.Range("F323").Value = 282.14

# Row 324
$ws.Range("A324").NumberFormat = $dateFormat
$ws.Range("A324").Value = 44264
$ws.Range("B324").Value = "Reko21"
$ws.Range("C324").Value = 2611
$ws.Range("D324").Value = "Reko Swish +46709927597"
$ws.Range("E324").Value = "'"
$ws.Range("F324").Value = 33.86

# Row 325
$ws.Range("A325").NumberFormat = $dateFormat
$ws.Range("A325").Value = 44264
$ws.Range("B325").Value = "Reko21"
$ws.Range("C325").Value = 1930
$ws.Range("D325").Value = "Reko Swish +46709927597"
$ws.Range("E325").Value = 316
$ws.Range("F325").Value = "'"

# Row 326
$ws.Range("A326").NumberFormat = $dateFormat
$ws.Range("A326").Value = 44264
$ws.Range("B326").Value = "'"
$ws.Range("C326").Value = 6400
$ws.Range("D326").Value = "FACEBK 7SAS3ZEZ62 K6885"
$ws.Range("E326").Value = 410
$ws.Range("F326").Value = "'"

# Row 327
$ws.Range("A327").NumberFormat = $dateFormat
$ws.Range("A327").Value = 44264
$ws.Range("B327").Value = "'"
$ws.Range("C327").Value = "'"
$ws.Range("D327").Value = "FACEBK 7SAS3ZEZ62 K6885"
$ws.Range("E327").Value = 0
$ws.Range("F327").Value = "'"

# Row 328
$ws.Range("A328").NumberFormat = $dateFormat
$ws.Range("A328").Value = 44264
$ws.Range("B328").Value = "'"
$ws.Range("C328").Value = 1930
$ws.Range("D328").Value = "FACEBK 7SAS3ZEZ62 K6885"
$ws.Range("E328").Value = "'"
$ws.Range("F328").Value = 410

# Row 329
$ws.Range("A329").NumberFormat = $dateFormat
$ws.Range("A329").Value = 44265
$ws.Range("B329").Value = "Reko22"
$ws.Range("C329").Value = 3011
$ws.Range("D329").Value = "Reko Swish +46708688090"
$ws.Range("E329").Value = "'"
$ws.Range("F329").Value = 230.36

# Row 330
$ws.Range("A330").NumberFormat = $dateFormat
$ws.Range("A330").Value = 44265
$ws.Range("B330").Value = "Reko22"
$ws.Range("C330").Value = 2611
$ws.Range("D330").Value = "Reko Swish +46708688090"
$ws.Range("E330").Value = "'"
$ws.Range("F330").Value = 27.64

# Row 331
$ws.Range("A331").NumberFormat = $dateFormat
$ws.Range("A331").Value = 44265
$ws.Range("B331").Value = "Reko22"
$ws.Range("C331").Value = 1930
$ws.Range("D331").Value = "Reko Swish +46708688090"
$ws.Range("E331").Value = 258
$ws.Range("F331").Value = "'"

# Row 332
$ws.Range("A332").NumberFormat = $dateFormat
$ws.Range("A332").Value = 44265
$ws.Range("B332").Value = "'3102253"
$ws.Range("C332").Value = 3011
$ws.Range("D332").Value = "Order 3102253 Card(Stripe)"
$ws.Range("E332").Value = "'"
$ws.Range("F332").Value = 431.25

# Row 333
$ws.Range("A333").NumberFormat = $dateFormat
$ws.Range("A333").Value = 44265
$ws.Range("B333").Value = "'3102253"
$ws.Range("C333").Value = 2611
$ws.Range("D333").Value = "Order 3102253 Card(Stripe)"
$ws.Range("E333").Value = "'"
$ws.Range("F333").Value = 51.75

# Row 334
$ws.Range("A334").NumberFormat = $dateFormat
$ws.Range("A334").Value = 44265
$ws.Range("B334").Value = "'3102253"
$ws.Range("C334").Value = 1930
$ws.Range("D334").Value = "Order 3102253 Card(Stripe)"
$ws.Range("E334").Value = 483
$ws.Range("F334").Value = "'"

# Row 335
$ws.Range("A335").NumberFormat = $dateFormat
$ws.Range("A335").Value = 44266
$ws.Range("B335").Value = "'8111258"
$ws.Range("C335").Value = 3011
$ws.Range("D335").Value = "Order 8111258 Card(Stripe)"
$ws.Range("E335").Value = "'"
$ws.Range("F335").Value = 1201.79

# Row 336
$ws.Range("A336").NumberFormat = $dateFormat
$ws.Range("A336").Value = 44266
$ws.Range("B336").Value = "'8111258"
$ws.Range("C336").Value = 2611
$ws.Range("D336").Value = "Order 8111258 Card(Stripe)"
$ws.Range("E336").Value = "'"
$ws.Range("F336").Value = 144.21

# Row 337
$ws.Range("A337").NumberFormat = $dateFormat
$ws.Range("A337").Value = 44266
$ws.Range("B337").Value = "'8111258"
$ws.Range("C337").Value = 1930
$ws.Range("D337").Value = "Order 8111258 Card(Stripe)"
$ws.Range("E337").Value = 1346
$ws.Range("F337").Value = "'"

# Row 338
$ws.Range("A338").NumberFormat = $dateFormat
$ws.Range("A338").Value = 44267
$ws.Range("B338").Value = "Reko23"
$ws.Range("C338").Value = 3011
$ws.Range("D338").Value = "Reko Swish +46702129177"
$ws.Range("E338").Value = "'"
$ws.Range("F338").Value = 345.54

# Row 339
$ws.Range("A339").NumberFormat = $dateFormat
$ws.Range("A339").Value = 44267
$ws.Range("B339").Value = "Reko23"
$ws.Range("C339").Value = 2611
$ws.Range("D339").Value = "Reko Swish +46702129177"
$ws.Range("E339").Value = "'"
$ws.Range("F339").Value = 41.46

# Row 340
$ws.Range("A340").NumberFormat = $dateFormat
$ws.Range("A340").Value = 44267
$ws.Range("B340").Value = "Reko23"
$ws.Range("C340").Value = 1930
$ws.Range("D340").Value = "Reko Swish +46702129177"
$ws.Range("E340").Value = 387
$ws.Range("F340").Value = "'"

# Row 341
$ws.Range("A341").NumberFormat = $dateFormat
$ws.Range("A341").Value = 44267
$ws.Range("B341").Value = "'"
$ws.Range("C341").Value = 4010
$ws.Range("D341").Value = "M&S RB BROMMA K0135"
$ws.Range("E341").Value = 437.67
$ws.Range("F341").Value = "'"

# Row 342
$ws.Range("A342").NumberFormat = $dateFormat
$ws.Range("A342").Value = 44267
$ws.Range("B342").Value = "'"
$ws.Range("C342").Value = 2645
$ws.Range("D342").Value = "M&S RB BROMMA K0135"
$ws.Range("E342").Value = 52.52
$ws.Range("F342").Value = "'"

# Row 343
$ws.Range("A343").NumberFormat = $dateFormat
$ws.Range("A343").Value = 44267
$ws.Range("B343").Value = "'"
$ws.Range("C343").Value = 1930
$ws.Range("D343").Value = "M&S RB BROMMA K0135"
$ws.Range("E343").Value = "'"
$ws.Range("F343").Value = 490.19

# Row 344
$ws.Range("A344").NumberFormat = $dateFormat
$ws.Range("A344").Value = 44268
$ws.Range("B344").Value = "'"
$ws.Range("C344").Value = 4010
$ws.Range("D344").Value = "KAHLS THE & KAFFEHANDE K6885"
$ws.Range("E344").Value = 269.64
$ws.Range("F344").Value = "'"

# Row 345
$ws.Range("A345").NumberFormat = $dateFormat
$ws.Range("A345").Value = 44268
$ws.Range("B345").Value = "'"
$ws.Range("C345").Value = 2645
$ws.Range("D345").Value = "KAHLS THE & KAFFEHANDE K6885"
$ws.Range("E345").Value = 32.36
$ws.Range("F345").Value = "'"

# Row 346
$ws.Range("A346").NumberFormat = $dateFormat
$ws.Range("A346").Value = 44268
$ws.Range("B346").Value = "'"
$ws.Range("C346").Value = 1930
$ws.Range("D346").Value = "KAHLS THE & KAFFEHANDE K6885"
$ws.Range("E346").Value = "'"
$ws.Range("F346").Value = 302

# Row 347
$ws.Range("A347").NumberFormat = $dateFormat
$ws.Range("A347").Value = 44269
$ws.Range("B347").Value = "'"
$ws.Range("C347").Value = 5010
$ws.Range("D347").Value = "Dec-MarKitchen"
$ws.Range("E347").Value = 15066
$ws.Range("F347").Value = "'"

# Row 348
$ws.Range("A348").NumberFormat = $dateFormat
$ws.Range("A348").Value = 44269
$ws.Range("B348").Value = "'"
$ws.Range("C348").Value = "'"
$ws.Range("D348").Value = "Dec-MarKitchen"
$ws.Range("E348").Value = 0
$ws.Range("F348").Value = "'"

# Row 349
$ws.Range("A349").NumberFormat = $dateFormat
$ws.Range("A349").Value = 44269
$ws.Range("B349").Value = "'"
$ws.Range("C349").Value = 1930
$ws.Range("D349").Value = "Dec-MarKitchen"
$ws.Range("E349").Value = "'"
$ws.Range("F349").Value = 15066

# Row 350
$ws.Range("A350").NumberFormat = $dateFormat
$ws.Range("A350").Value = 44269
$ws.Range("B350").Value = "'"
$ws.Range("C350").Value = 5670
$ws.Range("D350").Value = "ST1 V#LLINGBY K0135"
$ws.Range("E350").Value = 668.77
$ws.Range("F350").Value = "'"

# Row 351
$ws.Range("A351").NumberFormat = $dateFormat
$ws.Range("A351").Value = 44269
$ws.Range("B351").Value = "'"
$ws.Range("C351").Value = 2641
$ws.Range("D351").Value = "ST1 V#LLINGBY K0135"
$ws.Range("E351").Value = 167.19
$ws.Range("F351").Value = "'"

# Row 352
$ws.Range("A352").NumberFormat = $dateFormat
$ws.Range("A352").Value = 44269
$ws.Range("B352").Value = "'"
$ws.Range("C352").Value = 1930
$ws.Range("D352").Value = "ST1 V#LLINGBY K0135"
$ws.Range("E352").Value = "'"
$ws.Range("F352").Value = 835.96

# Row 353
$ws.Range("A353").NumberFormat = $dateFormat
$ws.Range("A353").Value = 44269
$ws.Range("B353").Value = "'"
$ws.Range("C353").Value = 6400
$ws.Range("D353").Value = "DECATHLON K0135"
$ws.Range("E353").Value = 319.2
$ws.Range("F353").Value = "'"

# Row 354
$ws.Range("A354").NumberFormat = $dateFormat
$ws.Range("A354").Value = 44269
$ws.Range("B354").Value = "'"
$ws.Range("C354").Value = 2641
$ws.Range("D354").Value = "DECATHLON K0135"
$ws.Range("E354").Value = 79.8
$ws.Range("F354").Value = "'"

# Row 355
$ws.Range("A355").NumberFormat = $dateFormat
$ws.Range("A355").Value = 44269
$ws.Range("B355").Value = "'"
$ws.Range("C355").Value = 1930
$ws.Range("D355").Value = "DECATHLON K0135"
$ws.Range("E355").Value = "'"
$ws.Range("F355").Value = 399
